# FDIOS-543: Removing personal information from the document templates
#
# The slide master and every slide layout carry an auto-updating
# "datetimeFigureOut" field placeholder (msoPlaceholderDate / PPT
# PlaceholderFormat.Type = 16) whose cached display text is "7/24/12".
# Refresh that cached text to "7/25/12" everywhere it appears.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $isDatePh = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
            if ($isDatePh) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "7/24/12") {
                    $tr.Text = "7/25/12"
                }
            }
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
